$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-07-24 Wednesday" "2024-07-25 Thursday"
Replace-Text "75×74=5550" "28×16=448"
Replace-Text "98×42=4116" "75×82=6150"
Replace-Text "95×23=2185" "97×91=8827"
Replace-Text "76×72=5472" "29×37=1073"
Replace-Text "36×84=3024" "52×70=3640"
Replace-Text "90×56=5040" "71×13=923"
Replace-Text "15×18=270" "73×77=5621"
Replace-Text "31×55=1705" "31×99=3069"
Replace-Text "37×68=2516" "14×64=896"
Replace-Text "17×69=1173" "77×33=2541"
Replace-Text "27×58=1566" "43×80=3440"
Replace-Text "30×26=780" "56×47=2632"
Replace-Text "34×97=3298" "69×29=2001"
Replace-Text "32×68=2176" "98×86=8428"
Replace-Text "18×52=936" "85×14=1190"
Replace-Text "63×51=3213" "27×76=2052"
Replace-Text "18×66=1188" "67×91=6097"
Replace-Text "40×84=3360" "28×75=2100"
Replace-Text "78×43=3354" "64×49=3136"
Replace-Text "73×15=1095" "84×65=5460"
Replace-Text "14×57=798" "98×95=9310"
Replace-Text "75×42=3150" "41×26=1066"
Replace-Text "94×22=2068" "33×78=2574"
Replace-Text "56×69=3864" "30×53=1590"
Replace-Text "98×21=2058" "62×11=682"

$d.Save()
